# HouseMarket.xlsx - add a second "Segment fullScan" block below the
# existing "Segment 1" table (rows 7-12), mirroring its layout/format,
# then leave the selection on J11 as the author did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: merged title banner, like row 1 but boxed on every cell ---
# Start from the existing boxed/shaded header cell (C1: bold, themed
# fill, thin box) so the new style reuses the same font/fill records,
# then square off the border to a full thin box on all four sides.
$ws.Range("C1").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)
$ws.Range("A7:C7").Borders.LineStyle = 1
$ws.Range("A7:C7").Borders.Weight = 2
$ws.Range("A7").Value = "Segment fullScan"
$ws.Range("A7:C7").Merge()

# --- Row 8: column headers, same look as row 2 ---
$ws.Range("A2:C2").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Attribute"
$ws.Range("B8").Value = "Min"
$ws.Range("C8").Value = "Max"

# --- Rows 9-12: data rows, same look as rows 3-5 ---
$ws.Range("A3:C3").Copy()
$ws.Range("A9:C12").PasteSpecial(-4122)

$ws.Range("A9").Value = "PriceUF"
$ws.Range("B9").Value = 1000
$ws.Range("C9").Value = 25000

$ws.Range("A10").Value = "MtTot"

$ws.Range("A11").Value = "Bdroom"

$ws.Range("A12").Value = "Region"
$ws.Range("B12").Value = "RM"

# Clear the clipboard marquee left by the last Copy()
$excel.CutCopyMode = 0

# Author finished with the selection sitting on J11
$ws.Range("J11").Select()
